$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3047.7632
$ws.Range("I15").Value = 3047.7632
$ws.Range("K15").Value = 9143.2896
$ws.Range("M15").Value = -8974.2896
$ws.Range("H111").Value = 1422.4584
$ws.Range("I111").Value = 1226.6842
$ws.Range("J111").Value = 2166.4
$ws.Range("K111").Value = 3680.0526
$ws.Range("L111").Value = 6499.200000000001
$ws.Range("M111").Value = -613.0526
$ws.Range("N111").Value = -12633.2
$ws.Range("H113").Value = 4680.5293
$ws.Range("J113").Value = 4763.3335
$ws.Range("L113").Value = 4763.3335
$ws.Range("N113").Value = -11271.3335
$ws.Range("H129").Value = 4167975.5
$ws.Range("J129").Value = 1295.909
$ws.Range("L129").Value = 3887.727
$ws.Range("N129").Value = -13887.727

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 20000
$ws.Range("I31").Value = 20000
$ws.Range("K31").Value = 20000
$ws.Range("M31").Value = -19706
$ws.Range("H32").Value = 3778.9192
$ws.Range("I32").Value = 3025.057
$ws.Range("J32").Value = 9809.817999999999
$ws.Range("K32").Value = 3025.057
$ws.Range("L32").Value = 9809.817999999999
$ws.Range("M32").Value = -2738.057
$ws.Range("N32").Value = -10383.818
$ws.Range("H63").Value = 2575
$ws.Range("I63").Value = 2475
$ws.Range("J63").Value = 2775
$ws.Range("K63").Value = 2475
$ws.Range("L63").Value = 2775
$ws.Range("M63").Value = -1789
$ws.Range("N63").Value = -4147
$ws.Range("H66").Value = 2575
$ws.Range("I66").Value = 2475
$ws.Range("J66").Value = 2775
$ws.Range("K66").Value = 12375
$ws.Range("L66").Value = 13875
$ws.Range("M66").Value = -8943
$ws.Range("N66").Value = -20739
$ws.Range("H74").Value = 786.8333
$ws.Range("I74").Value = 731.9048
$ws.Range("J74").Value = 1171.3334
$ws.Range("K74").Value = 731.9048
$ws.Range("L74").Value = 1171.3334
$ws.Range("M74").Value = 142.0952
$ws.Range("N74").Value = -2919.3334
$ws.Range("H77").Value = 786.8333
$ws.Range("I77").Value = 731.9048
$ws.Range("J77").Value = 1171.3334
$ws.Range("K77").Value = 3659.524
$ws.Range("L77").Value = 5856.666999999999
$ws.Range("M77").Value = 708.4759999999997
$ws.Range("N77").Value = -14592.667
$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 1866.6666
$ws.Range("J88").Value = 2200
$ws.Range("K88").Value = 1866.6666
$ws.Range("L88").Value = 2200
$ws.Range("M88").Value = -1460.6666
$ws.Range("N88").Value = -3012
$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 1866.6666
$ws.Range("J91").Value = 2200
$ws.Range("K91").Value = 1866.6666
$ws.Range("L91").Value = 2200
$ws.Range("M91").Value = -462.6666
$ws.Range("N91").Value = -5008
$ws.Range("H102").Value = 3359.9333
$ws.Range("I102").Value = 3299.8572
$ws.Range("J102").Value = 3412.5
$ws.Range("K102").Value = 3299.8572
$ws.Range("L102").Value = 3412.5
$ws.Range("M102").Value = -1677.8572
$ws.Range("N102").Value = -6656.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2160
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H58").Value = 11907629
$ws.Range("I58").Value = 1930.96
$ws.Range("J58").Value = 29416008
$ws.Range("K58").Value = 1930.96
$ws.Range("L58").Value = 29416008
$ws.Range("M58").Value = -1727.96
$ws.Range("N58").Value = -29416414
$ws.Range("H74").Value = 20866.4
$ws.Range("J74").Value = 20866.4
$ws.Range("L74").Value = 20866.4
$ws.Range("N74").Value = -22614.4
$ws.Range("H77").Value = 20866.4
$ws.Range("J77").Value = 20866.4
$ws.Range("L77").Value = 62599.2
$ws.Range("N77").Value = -71335.20000000001
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H105").Value = 2629.9375
$ws.Range("J105").Value = 2150
$ws.Range("L105").Value = 2150
$ws.Range("N105").Value = -5644
$ws.Range("H107").Value = 1848.6316
$ws.Range("I107").Value = 446
$ws.Range("K107").Value = 446
$ws.Range("M107").Value = 1474
$ws.Range("H113").Value = 2160
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 2511.1667
$ws.Range("I122").Value = 2376.2856
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 7128.8568
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -4678.8568
$ws.Range("N122").Value = -13000
$ws.Range("H132").Value = 2553.8647
$ws.Range("J132").Value = 4156
$ws.Range("L132").Value = 12468
$ws.Range("N132").Value = -17528
$ws.Range("H136").Value = 11907629
$ws.Range("I136").Value = 1930.96
$ws.Range("J136").Value = 29416008
$ws.Range("K136").Value = 5792.88
$ws.Range("L136").Value = 88248024
$ws.Range("M136").Value = -3242.88
$ws.Range("N136").Value = -88253124

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1115
$ws.Range("I86").Value = 547.5
$ws.Range("K86").Value = 1642.5
$ws.Range("M86").Value = -456.5
$ws.Range("H89").Value = 1115
$ws.Range("I89").Value = 547.5
$ws.Range("K89").Value = 4927.5
$ws.Range("M89").Value = 1000.5
$ws.Range("H92").Value = 3185.7144
$ws.Range("J92").Value = 3633.3333
$ws.Range("L92").Value = 10899.9999
$ws.Range("N92").Value = -13395.9999
$ws.Range("H122").Value = 1418.125
$ws.Range("I122").Value = 801.3333
$ws.Range("J122").Value = 1560.4615
$ws.Range("K122").Value = 7211.9997
$ws.Range("L122").Value = 14044.1535
$ws.Range("M122").Value = -4761.9997
$ws.Range("N122").Value = -18944.1535
$ws.Range("H132").Value = 2437.0908
$ws.Range("I132").Value = 1750.5
$ws.Range("J132").Value = 3261
$ws.Range("K132").Value = 15754.5
$ws.Range("L132").Value = 29349
$ws.Range("M132").Value = -13224.5
$ws.Range("N132").Value = -34409

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 60003
$ws.Range("J19").Value = 63337.332
$ws.Range("L19").Value = 63337.332
$ws.Range("N19").Value = -63913.332
$ws.Range("H97").Value = 2068.25
$ws.Range("I97").Value = 1331.9
$ws.Range("J97").Value = 5750
$ws.Range("K97").Value = 1331.9
$ws.Range("L97").Value = 5750
$ws.Range("M97").Value = -835.9000000000001
$ws.Range("N97").Value = -6742

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 90912630
$ws.Range("I61").Value = 125003260
$ws.Range("J61").Value = 4268.3335
$ws.Range("K61").Value = 125003260
$ws.Range("L61").Value = 4268.3335
$ws.Range("M61").Value = -125003058
$ws.Range("N61").Value = -4672.3335
$ws.Range("H113").Value = 90912630
$ws.Range("I113").Value = 125003260
$ws.Range("J113").Value = 4268.3335
$ws.Range("K113").Value = 125003260
$ws.Range("L113").Value = 4268.3335
$ws.Range("M113").Value = -125001090
$ws.Range("N113").Value = -8608.333500000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 478481.53
$ws.Range("I122").Value = 626694.5
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 1880083.5
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -1877633.5
$ws.Range("N122").Value = -17500
$ws.Range("H132").Value = 208115.45
$ws.Range("J132").Value = 12745.728
$ws.Range("L132").Value = 38237.18399999999
$ws.Range("N132").Value = -43297.18399999999
